$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextValue "D2" "328.69"
Set-TextValue "E2" "-0.64%"
Set-TextValue "D3" "43.88"
Set-TextValue "E3" "5.26%"
Set-TextValue "D4" "5.425"
Set-TextValue "E4" "-4.55%"
Set-TextValue "D5" "0.08106"
Set-TextValue "E5" "-2.91%"
Set-TextValue "D6" "8.695"
Set-TextValue "E6" "-1.23%"
Set-TextValue "D7" "4.303"
Set-TextValue "E7" "-3.65%"
Set-TextValue "D8" "1.900"
Set-TextValue "E8" "-4.55%"
Set-TextValue "D10" "0.9436"
Set-TextValue "E10" "1.87%"
Set-TextValue "D11" "0.1189"
Set-TextValue "E11" "-8.10%"
Set-TextValue "E12" "-4.63%"
Set-TextValue "D13" "0.09621"
Set-TextValue "E13" "0.12%"
Set-TextValue "D14" "0.04220"
Set-TextValue "E14" "9.10%"
Set-TextValue "D15" "0.1070"
Set-TextValue "E15" "0.93%"
Set-TextValue "D16" "0.001275"
Set-TextValue "E16" "-2.16%"
Set-TextValue "D17" "0.005969"
Set-TextValue "E17" "-2.26%"
Set-TextValue "E18" "3.72%"
Set-TextValue "D20" "8.600"
Set-TextValue "E20" "0.89%"
Set-TextValue "D21" "0.1360"
Set-TextValue "E21" "-0.20%"
Set-TextValue "D22" "0.2605"
Set-TextValue "E22" "4.86%"
Set-TextValue "D23" "0.04391"
Set-TextValue "E23" "-0.43%"
Set-TextValue "D24" "0.001242"
Set-TextValue "E24" "-2.86%"
Set-TextValue "D25" "0.004328"
Set-TextValue "E25" "-1.07%"
Set-TextValue "E26" "1.44%"
Set-TextValue "E27" "31.90%"
Set-TextValue "D39" "0.02701"
Set-TextValue "E39" "-4.94%"
Set-TextValue "D40" "0.05543"
Set-TextValue "E40" "0.42%"
Set-TextValue "D41" "0.007803"
Set-TextValue "E41" "-1.82%"
Set-TextValue "D42" "0.009763"
Set-TextValue "E42" "4.92%"
Set-TextValue "D43" "0.1398"
Set-TextValue "E43" "-2.58%"
Set-TextValue "E44" "-0.73%"
Set-TextValue "D45" "0.009638"
Set-TextValue "E45" "-12.81%"
Set-TextValue "D46" "0.00007108"
Set-TextValue "E46" "1.83%"
Set-TextValue "E47" "0.66%"
Set-TextValue "D48" "0.003476"
Set-TextValue "E48" "0.75%"
Set-TextValue "E49" "0.32%"
Set-TextValue "E50" "0.66%"
Set-TextValue "E51" "0.66%"
